# Natmi following Dr Hou advice
# Update LR-pairs_lrc2p Adm-Calcr results: recompute rows 2-4 and add new rows 5-9
# for the additional "M1" sending/target cluster.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: ECs -> ECs
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Adm"
$ws.Range("C2").Value = "Calcr"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 12.600084
$ws.Range("H2").Value = 37.800252
$ws.Range("I2").Value = 0.4109984227877579
$ws.Range("J2").Value = 0.4109984227877579
$ws.Range("K2").Value = 1
$ws.Range("L2").Value = 0.3333333333333333
$ws.Range("M2").Value = 0.096149
$ws.Range("N2").Value = 0.288447
$ws.Range("O2").Value = 0.2085078145161115
$ws.Range("P2").Value = 0.2085078145161116
$ws.Range("Q2").Value = 1.211485476516
$ws.Range("R2").Value = 10.903369288644
$ws.Range("S2").Value = 0.08569638290504421
$ws.Range("T2").Value = 0.08569638290504422

# Row 3: ECs -> sCs
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Adm"
$ws.Range("C3").Value = "Calcr"
$ws.Range("D3").Value = "sCs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 12.600084
$ws.Range("H3").Value = 37.800252
$ws.Range("I3").Value = 0.4109984227877579
$ws.Range("J3").Value = 0.4109984227877579
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 0.36498
$ws.Range("N3").Value = 1.09494
$ws.Range("O3").Value = 0.7914921854838884
$ws.Range("P3").Value = 0.7914921854838886
$ws.Range("Q3").Value = 4.598778658320001
$ws.Range("R3").Value = 41.38900792488
$ws.Range("S3").Value = 0.3253020398827137
$ws.Range("T3").Value = 0.3253020398827137

# Row 4: FAPs -> ECs
$ws.Range("A4").Value = "FAPs"
$ws.Range("B4").Value = "Adm"
$ws.Range("C4").Value = "Calcr"
$ws.Range("D4").Value = "ECs"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 17.59249933333333
$ws.Range("H4").Value = 52.77749799999999
$ws.Range("I4").Value = 0.573844545710543
$ws.Range("J4").Value = 0.573844545710543
$ws.Range("K4").Value = 1
$ws.Range("L4").Value = 0.3333333333333333
$ws.Range("M4").Value = 0.096149
$ws.Range("N4").Value = 0.288447
$ws.Range("O4").Value = 0.2085078145161115
$ws.Range("P4").Value = 0.2085078145161116
$ws.Range("Q4").Value = 1.691501218400667
$ws.Range("R4").Value = 15.223510965606
$ws.Range("S4").Value = 0.1196510720980962
$ws.Range("T4").Value = 0.1196510720980962

# Row 5: FAPs -> sCs
$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Adm"
$ws.Range("C5").Value = "Calcr"
$ws.Range("D5").Value = "sCs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 17.59249933333333
$ws.Range("H5").Value = 52.77749799999999
$ws.Range("I5").Value = 0.573844545710543
$ws.Range("J5").Value = 0.573844545710543
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 0.36498
$ws.Range("N5").Value = 1.09494
$ws.Range("O5").Value = 0.7914921854838884
$ws.Range("P5").Value = 0.7914921854838886
$ws.Range("Q5").Value = 6.42091040668
$ws.Range("R5").Value = 57.78819366011999
$ws.Range("S5").Value = 0.4541934736124468
$ws.Range("T5").Value = 0.4541934736124469

# Row 6: M1 -> ECs
$ws.Range("A6").Value = "M1"
$ws.Range("B6").Value = "Adm"
$ws.Range("C6").Value = "Calcr"
$ws.Range("D6").Value = "ECs"
$ws.Range("E6").Value = 1
$ws.Range("F6").Value = 0.3333333333333333
$ws.Range("G6").Value = 0.02822
$ws.Range("H6").Value = 0.08466
$ws.Range("I6").Value = 0.0009204998546891057
$ws.Range("J6").Value = 0.0009204998546891058
$ws.Range("K6").Value = 1
$ws.Range("L6").Value = 0.3333333333333333
$ws.Range("M6").Value = 0.096149
$ws.Range("N6").Value = 0.288447
$ws.Range("O6").Value = 0.2085078145161115
$ws.Range("P6").Value = 0.2085078145161116
$ws.Range("Q6").Value = 0.00271332478
$ws.Range("R6").Value = 0.02441992302
$ws.Range("S6").Value = 0.0001919314129636237
$ws.Range("T6").Value = 0.0001919314129636237

# Row 7: M1 -> sCs
$ws.Range("A7").Value = "M1"
$ws.Range("B7").Value = "Adm"
$ws.Range("C7").Value = "Calcr"
$ws.Range("D7").Value = "sCs"
$ws.Range("E7").Value = 1
$ws.Range("F7").Value = 0.3333333333333333
$ws.Range("G7").Value = 0.02822
$ws.Range("H7").Value = 0.08466
$ws.Range("I7").Value = 0.0009204998546891057
$ws.Range("J7").Value = 0.0009204998546891058
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 0.36498
$ws.Range("N7").Value = 1.09494
$ws.Range("O7").Value = 0.7914921854838884
$ws.Range("P7").Value = 0.7914921854838886
$ws.Range("Q7").Value = 0.0102997356
$ws.Range("R7").Value = 0.0926976204
$ws.Range("S7").Value = 0.000728568441725482
$ws.Range("T7").Value = 0.0007285684417254822

# Row 8: sCs -> ECs
$ws.Range("A8").Value = "sCs"
$ws.Range("B8").Value = "Adm"
$ws.Range("C8").Value = "Calcr"
$ws.Range("D8").Value = "ECs"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 0.436453
$ws.Range("H8").Value = 1.309359
$ws.Range("I8").Value = 0.01423653164701007
$ws.Range("J8").Value = 0.01423653164701007
$ws.Range("K8").Value = 1
$ws.Range("L8").Value = 0.3333333333333333
$ws.Range("M8").Value = 0.096149
$ws.Range("N8").Value = 0.288447
$ws.Range("O8").Value = 0.2085078145161115
$ws.Range("P8").Value = 0.2085078145161116
$ws.Range("Q8").Value = 0.041964519497
$ws.Range("R8").Value = 0.377680675473
$ws.Range("S8").Value = 0.002968428100007528
$ws.Range("T8").Value = 0.002968428100007529

# Row 9: sCs -> sCs
$ws.Range("A9").Value = "sCs"
$ws.Range("B9").Value = "Adm"
$ws.Range("C9").Value = "Calcr"
$ws.Range("D9").Value = "sCs"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 0.436453
$ws.Range("H9").Value = 1.309359
$ws.Range("I9").Value = 0.01423653164701007
$ws.Range("J9").Value = 0.01423653164701007
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 0.36498
$ws.Range("N9").Value = 1.09494
$ws.Range("O9").Value = 0.7914921854838884
$ws.Range("P9").Value = 0.7914921854838886
$ws.Range("Q9").Value = 0.15929661594
$ws.Range("R9").Value = 1.43366954346
$ws.Range("S9").Value = 0.01126810354700254
$ws.Range("T9").Value = 0.01126810354700255

